$wb = $excel.ActiveWorkbook

$wsSearch = $wb.Worksheets.Item("SearchData")
$wsTest   = $wb.Worksheets.Item("Testing")

# --- SearchData sheet -------------------------------------------------
# Bump the "After Days" search parameter so the date filter skips the
# unwanted dates from last month (21 -> 104 days out).
$wsSearch.Range("C2").Value = 104

# Park the cursor on the changed cell (mirrors what the author did in
# Excel) without disturbing which sheet tab is active overall.
$excel.Goto($wsSearch.Range("C2"))

# --- Testing sheet ------------------------------------------------------
# The "Expected" date for the Date test case needs to reflect the new
# target date produced by the updated After-Days value.
$wsTest.Range("C6").Value = "05-Oct-25 Sun"

# Extend the hidden "no scroll past this" formatting block below the
# visible data by a further six rows (25 -> 31) to keep stray dates from
# last month out of view, matching the larger used range.
$wsTest.Rows("26:31").Hidden = $true

# Restore "Testing" as the active/visible tab.
$wsTest.Activate()

$wb.Save()
